$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blatt1")

# Collapse the duplicated "Besuchter Link" / "Link" cell styles that had
# accumulated (Excel keeps only one set per workbook); trim back down to
# 5 of each plus "Standard".
for ($i = 32; $i -ge 22; $i--) {
  $wb.Styles.Item($i).Delete()
}
for ($i = 16; $i -ge 6; $i--) {
  $wb.Styles.Item($i).Delete()
}
for ($i = 1; $i -le $wb.Styles.Count; $i++) {
  $s = $wb.Styles.Item($i)
  if ($s.Name -eq "Link" -or $s.Name -eq "Besuchter Link") {
    $s.Hidden = $true
  }
}

# Update data values
$ws.Range("L4").Value = 378.89699999999999
$ws.Range("L5").Value = 363.97800000000001
$ws.Range("L6").Value = 374.30099999999999

$ws.Range("G12").Value = 553.23699999999997
$ws.Range("R12").Value = 723.73299999999995

$ws.Range("G13").Value = 533.92399999999998
$ws.Range("R13").ClearContents()

$ws.Range("G14").Value = 537.89300000000003

# Remove the now-stale ratio formulas (data removed alongside underlying timing changes)
$ws.Range("O4").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("U12").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("U13").ClearContents()

# Update selection to match final author state
$ws.Range("R12").Select()

# Turn the sheet ruler back on (was previously turned off)
$excel.ActiveWindow.DisplayRuler = $true
